$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header text fixes.
#    "Из них негодно" (shared string formerly used by H1) is renamed to
#    "Брак", and the two rightmost data columns swap their header meaning:
#    H1 becomes "Цена" and I1 becomes "Брак".
# ---------------------------------------------------------------------------
$ws.Range("H1").Value = "Цена"
$ws.Range("I1").Value = "Брак"

# ---------------------------------------------------------------------------
# 2. Header row formatting: bold font, yellow fill, thin black border all
#    around, centered horizontally and vertically.
#    The format is assembled once on a scratch cell and then applied to the
#    whole header row (A1:J1) in a single PasteSpecial so the workbook ends
#    up with one clean combined style instead of many intermediate ones.
# ---------------------------------------------------------------------------
$scratch = $ws.Range("Z1")
$scratch.Style = "Normal"
$scratch.Font.Bold = $true
$scratch.Interior.Color = 65535
$scratch.Borders.LineStyle = 1
$scratch.Borders.Weight = 2
$scratch.HorizontalAlignment = -4108
$scratch.VerticalAlignment = -4108

$scratch.Copy() | Out-Null
$headerRng = $ws.Range("A1:J1")
$headerRng.PasteSpecial(-4122) | Out-Null
$scratch.Clear() | Out-Null

# ---------------------------------------------------------------------------
# 3. Header row is taller now.
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 30

# ---------------------------------------------------------------------------
# 4. Page is set up for portrait printing.
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 5. Leave the selection where the author left it.
# ---------------------------------------------------------------------------
$ws.Range("H23").Select() | Out-Null
